# Apply the cryptos.xlsx data refresh (GitHub Actions scheduled update).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "43.020.55"
$ws.Range("E2").Value = "  -0.14%  "
$ws.Range("D3").Value = "2.307.95"
$ws.Range("E3").Value = "  +0.08%  "
$ws.Range("E4").Value = "  +0.02%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "300.11"
$ws.Range("E5").Value = "  -0.51%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "98.00"
$ws.Range("E6").Value = "  -0.80%  "
$ws.Range("E7").Value = "  -2.14%  "
$ws.Range("E8").Value = "  +0.03%  "
$ws.Range("E9").Value = "  -2.98%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "36.06"
$ws.Range("E10").Value = "  +0.78%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0791"
$ws.Range("E11").Value = "  -0.24%  "
$ws.Range("E12").Value = "  +1.85%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "18.08"
$ws.Range("E13").Value = "  +1.26%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "6.79"
$ws.Range("E14").Value = "  -1.75%  "
$ws.Range("D15").Value = "2.668.16"
$ws.Range("E15").Value = "  +0.13%  "
$ws.Range("D16").Value = "2.306.01"
$ws.Range("E16").Value = "  +1.86%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.782"
$ws.Range("E17").Value = "  -0.97%  "
$ws.Range("D18").Value = "42.957.43"
$ws.Range("E18").Value = "  -0.04%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "12.78"
$ws.Range("E19").Value = "  -5.17%  "
$ws.Range("D20").Value = "0.0₃0904"
$ws.Range("E20").Value = "  -0.65%  "
$ws.Range("E21").Value = "  -1.99%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "67.92"
$ws.Range("E22").Value = "  -0.64%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "237.25"
$ws.Range("E23").Value = "  -0.94%  "
$ws.Range("E24").Value = "  -1.25%  "
$ws.Range("E26").Value = "  +0.09%  "
$ws.Range("E27").Value = "  -0.54%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "25.40"
$ws.Range("E28").Value = "  +2.33%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "165.15"
$ws.Range("E29").Value = "  -1.50%  "
$ws.Range("B30").Value = "Toncoin"
$ws.Range("C30").Value = "https://coinranking.com/coin/67YlI0K1b+toncoin-ton"
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "2.03"
$ws.Range("E30").Value = "  -0.34%  "
$ws.Range("B31").Value = "Cosmos"
$ws.Range("C31").Value = "https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom"
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "9.07"
$ws.Range("E31").Value = "  -0.73%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "33.30"
$ws.Range("E32").Value = "  -0.15%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "4.92"
$ws.Range("E33").Value = "  +1.72%  "
$ws.Range("E34").Value = "  +0.08%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "5.03"
$ws.Range("E35").Value = "  -3.98%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "17.02"
$ws.Range("E36").Value = "  -5.96%  "
$ws.Range("E37").Value = "  -1.04%  "
$ws.Range("E38").Value = "  -0.12%  "
$ws.Range("E39").Value = "  -0.78%  "
$ws.Range("E40").Value = "  -1.59%  "
$ws.Range("E41").Value = "  -1.78%  "
$ws.Range("E42").Value = "  -0.60%  "
$ws.Range("D43").Value = "2.019.48"
$ws.Range("E43").Value = "  +0.75%  "
$ws.Range("E44").Value = "  -2.13%  "
$ws.Range("E45").Value = "  +0.77%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "10.09"
$ws.Range("E46").Value = "  -0.11%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "17.44"
$ws.Range("E47").Value = "  -0.04%  "
$ws.Range("E48").Value = "  -1.42%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "2.93"
$ws.Range("E49").Value = "  -0.08%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "53.88"
$ws.Range("E50").Value = "  -1.29%  "
$ws.Range("D51").Value = "2.536.29"
$ws.Range("E51").Value = "  +0.28%  "
